$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry in the version history table
$ws.Range("A17").Value = 45751
$ws.Range("A5").Copy()
$ws.Range("A17").PasteSpecial(-4122)

$ws.Range("B17").Value = 0.79166666666666663
$ws.Range("B6").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("C17").Value = "Futconnect0404 1900"
$ws.Range("D17").Value = "Ajustes nos filtros de datas e no apelido e email ao cadastro um novo sócio."

# Move the selection to the next empty row, matching the post-edit state
$ws.Range("D18").Select()
